$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.117959499359131
$ws.Range("B1").Value = 1.012356162071228
$ws.Range("C1").Value = 6.664680480957031
$ws.Range("D1").Value = 2.025692462921143
$ws.Range("E1").Value = 1.125463128089905
